$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# New column N: "Levy" header + values
$ws2.Range("N3").Value = "Levy"
$ws2.Range("N4").Value = 0
$ws2.Range("N5").Value = 30
$ws2.Range("N6").Value = 35

# New column F: initial VLOOKUP formulas
$ws2.Range("F4").Formula = '=VLOOKUP(E4,$L$4:$N$6,2)'
$ws2.Range("F5").Formula = '=VLOOKUP(E5,$L$4:$N$6,2)'
$ws2.Range("F6").Formula = '=VLOOKUP(E6,$L$4:$N$6,2)'

# Activate Sheet2 and set selection to E12 (matches final workbook view state)
$ws2.Activate() | Out-Null
$ws2.Range("E12").Select() | Out-Null
